$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new gyroscope sample was recorded and prepended to the dataset; the
# rest of the samples shift down one row and the single oldest sample
# (previously the last row, 22) falls off the end.

# Capture the existing x/y/z samples (rows 2-21) before overwriting them.
$existing = @{}
for ($r = 2; $r -le 21; $r++) {
    $existing[$r] = @($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 2).Value2, $ws.Cells.Item($r, 3).Value2)
}

# Row 2 gets the newly recorded sample.
$ws.Cells.Item(2, 1).Value = -0.2335032373666763
$ws.Cells.Item(2, 2).Value = -0.1345430761575698
$ws.Cells.Item(2, 3).Value = 0.1078177168965339

# Rows 3-21 get what used to be in rows 2-20 (the old row 21 sample is
# dropped along with the old row 22 sample).
for ($r = 21; $r -ge 3; $r--) {
    $src = $existing[$r - 1]
    $ws.Cells.Item($r, 1).Value = $src[0]
    $ws.Cells.Item($r, 2).Value = $src[1]
    $ws.Cells.Item($r, 3).Value = $src[2]
}

# Remove the now-stale trailing row so the dataset stays at 20 rows
# (A2:C21).
$ws.Range("A22:C22").EntireRow.Delete()
